$d = $word.ActiveDocument
$d.Content.Find.Execute("322×6=", $true, $false, $false, $false, $false, $true, 1, $false, "957×6=", 2) | Out-Null
$d.Content.Find.Execute("213×6=", $true, $false, $false, $false, $false, $true, 1, $false, "555×7=", 2) | Out-Null
$d.Content.Find.Execute("353×2=", $true, $false, $false, $false, $false, $true, 1, $false, "355×4=", 2) | Out-Null
$d.Content.Find.Execute("359×6=", $true, $false, $false, $false, $false, $true, 1, $false, "638×6=", 2) | Out-Null
$d.Content.Find.Execute("572×5=", $true, $false, $false, $false, $false, $true, 1, $false, "129×9=", 2) | Out-Null
$d.Content.Find.Execute("170×6=", $true, $false, $false, $false, $false, $true, 1, $false, "399×7=", 2) | Out-Null
$d.Content.Find.Execute("475×2=", $true, $false, $false, $false, $false, $true, 1, $false, "523×6=", 2) | Out-Null
$d.Content.Find.Execute("251×7=", $true, $false, $false, $false, $false, $true, 1, $false, "686×4=", 2) | Out-Null
$d.Content.Find.Execute("849×2=", $true, $false, $false, $false, $false, $true, 1, $false, "773×8=", 2) | Out-Null
$d.Content.Find.Execute("963×6=", $true, $false, $false, $false, $false, $true, 1, $false, "657×7=", 2) | Out-Null
$d.Content.Find.Execute("286×9=", $true, $false, $false, $false, $false, $true, 1, $false, "369×8=", 2) | Out-Null
$d.Content.Find.Execute("914×2=", $true, $false, $false, $false, $false, $true, 1, $false, "440×5=", 2) | Out-Null
$d.Content.Find.Execute("805×5=", $true, $false, $false, $false, $false, $true, 1, $false, "204×4=", 2) | Out-Null
$d.Content.Find.Execute("780×3=", $true, $false, $false, $false, $false, $true, 1, $false, "180×7=", 2) | Out-Null
$d.Content.Find.Execute("540×3=", $true, $false, $false, $false, $false, $true, 1, $false, "752×5=", 2) | Out-Null
$d.Content.Find.Execute("512×3=", $true, $false, $false, $false, $false, $true, 1, $false, "889×5=", 2) | Out-Null
$d.Content.Find.Execute("785×9=", $true, $false, $false, $false, $false, $true, 1, $false, "893×5=", 2) | Out-Null
$d.Content.Find.Execute("842×4=", $true, $false, $false, $false, $false, $true, 1, $false, "684×3=", 2) | Out-Null
$d.Content.Find.Execute("260×5=", $true, $false, $false, $false, $false, $true, 1, $false, "952×7=", 2) | Out-Null
$d.Content.Find.Execute("856×3=", $true, $false, $false, $false, $false, $true, 1, $false, "242×3=", 2) | Out-Null
$d.Content.Find.Execute("857×3=", $true, $false, $false, $false, $false, $true, 1, $false, "419×8=", 2) | Out-Null
$d.Content.Find.Execute("583×5=", $true, $false, $false, $false, $false, $true, 1, $false, "949×4=", 2) | Out-Null
$d.Content.Find.Execute("759×4=", $true, $false, $false, $false, $false, $true, 1, $false, "178×6=", 2) | Out-Null
$d.Content.Find.Execute("249×5=", $true, $false, $false, $false, $false, $true, 1, $false, "522×5=", 2) | Out-Null
$d.Content.Find.Execute("630×6=", $true, $false, $false, $false, $false, $true, 1, $false, "262×4=", 2) | Out-Null
$d.Save()
